$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Copy the header style from C1 (Score) onto the new D1 header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "country_code"

# New "country_code" column data, mirroring the Score column values.
$ws.Range("D2").Value = 90
$ws.Range("D3").Value = 99
$ws.Range("D4").Value = 50
